# Apply cryptos list update (Wed Feb 14 15:06:27 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Worksheet, $Address, $Val) {
    $rng = $Worksheet.Range($Address)
    $rng.NumberFormat = "@"
    $rng.Value = $Val
    $rng.Style = "Normal"
}

Set-TextValue $ws "D2" "51.886.05"
Set-TextValue $ws "E2" "  +6.38%  "

Set-TextValue $ws "D3" "2.759.43"
Set-TextValue $ws "E3" "  +4.70%  "

Set-TextValue $ws "D4" "1.00"
Set-TextValue $ws "E4" "  -0.03%  "

Set-TextValue $ws "D5" "117.73"
Set-TextValue $ws "E5" "  +7.05%  "

Set-TextValue $ws "D6" "332.28"
Set-TextValue $ws "E6" "  +3.28%  "

Set-TextValue $ws "E7" "  +2.97%  "

Set-TextValue $ws "D8" "1.00"
Set-TextValue $ws "E8" "  +0.10%  "

Set-TextValue $ws "D9" "0.577"
Set-TextValue $ws "E9" "  +7.13%  "

Set-TextValue $ws "D10" "41.81"
Set-TextValue $ws "E10" "  +5.99%  "

Set-TextValue $ws "D11" "20.03"
Set-TextValue $ws "E11" "  +1.58%  "

Set-TextValue $ws "D12" "0.0831"
Set-TextValue $ws "E12" "  +2.82%  "

Set-TextValue $ws "E13" "  +3.05%  "

Set-TextValue $ws "D14" "7.63"
Set-TextValue $ws "E14" "  +5.89%  "

Set-TextValue $ws "D15" "3.192.50"
Set-TextValue $ws "E15" "  +4.88%  "

Set-TextValue $ws "D16" "2.767.56"
Set-TextValue $ws "E16" "  +4.64%  "

Set-TextValue $ws "D17" "0.884"
Set-TextValue $ws "E17" "  +3.08%  "

Set-TextValue $ws "D18" "51.767.65"
Set-TextValue $ws "E18" "  +6.13%  "

Set-TextValue $ws "D19" "13.63"
Set-TextValue $ws "E19" "  +6.52%  "

Set-TextValue $ws "B20" "ImmutableX"
Set-TextValue $ws "C20" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws "D20" "3.00"
Set-TextValue $ws "E20" "  +3.56%  "

Set-TextValue $ws "B21" "Uniswap"
Set-TextValue $ws "C21" "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue $ws "D21" "6.86"
Set-TextValue $ws "E21" "  +3.12%  "

Set-TextValue $ws "D22" "0.0₃0966"
Set-TextValue $ws "E22" "  +2.83%  "

Set-TextValue $ws "D23" "279.20"
Set-TextValue $ws "E23" "  +3.60%  "

Set-TextValue $ws "D24" "69.85"
Set-TextValue $ws "E24" "  +0.27%  "

Set-TextValue $ws "D25" "2.64"
Set-TextValue $ws "E25" "  +4.43%  "

Set-TextValue $ws "D26" "26.87"
Set-TextValue $ws "E26" "  +2.85%  "

Set-TextValue $ws "D27" "4.15"
Set-TextValue $ws "E27" "  +0.44%  "

Set-TextValue $ws "D28" "0.999"
Set-TextValue $ws "E28" "  +0.01%  "

Set-TextValue $ws "D29" "10.26"
Set-TextValue $ws "E29" "  +2.10%  "

Set-TextValue $ws "E30" "  +0.23%  "

Set-TextValue $ws "E31" "  +3.07%  "

Set-TextValue $ws "D32" "35.28"
Set-TextValue $ws "E32" "  +0.56%  "

Set-TextValue $ws "D33" "50.45"
Set-TextValue $ws "E33" "  +2.40%  "

Set-TextValue $ws "D34" "5.59"
Set-TextValue $ws "E34" "  +3.46%  "

Set-TextValue $ws "E35" "  +3.77%  "

Set-TextValue $ws "D36" "19.12"
Set-TextValue $ws "E36" "  -0.29%  "

Set-TextValue $ws "E37" "  -0.18%  "

Set-TextValue $ws "E38" "  +4.22%  "

Set-TextValue $ws "D39" "5.00"
Set-TextValue $ws "E39" "  +1.56%  "

Set-TextValue $ws "D40" "3.22"
Set-TextValue $ws "E40" "  +3.04%  "

Set-TextValue $ws "D41" "131.77"
Set-TextValue $ws "E41" "  +5.61%  "

Set-TextValue $ws "D42" "23.24"
Set-TextValue $ws "E42" "  +2.41%  "

Set-TextValue $ws "E43" "  +10.86%  "

Set-TextValue $ws "E44" "  +2.82%  "

Set-TextValue $ws "E45" "  +5.73%  "

Set-TextValue $ws "D46" "2.39"
Set-TextValue $ws "E46" "  +14.64%  "

Set-TextValue $ws "D47" "2.112.79"
Set-TextValue $ws "E47" "  +2.18%  "

Set-TextValue $ws "D48" "3.34"
Set-TextValue $ws "E48" "  +4.35%  "

Set-TextValue $ws "D49" "2.25"
Set-TextValue $ws "E49" "  +3.00%  "

Set-TextValue $ws "D50" "5.60"
Set-TextValue $ws "E50" "  +8.56%  "

Set-TextValue $ws "D51" "8.98"
Set-TextValue $ws "E51" "  +0.60%  "
